$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 5 email address first (new shared string insertion order)
$ws.Range("B5").Value = "sj@g.com"

# Add header for new "remarks" column
$ws.Range("C1").Value = "remarks"

# Add remarks values for existing test rows
$ws.Range("C2").Value = "Invalid"
$ws.Range("C3").Value = "Invalid"
$ws.Range("C4").Value = "Invalid"

# Mark row 5 as valid
$ws.Range("C5").Value = "Valid"

# Update selection to match target state
$ws.Range("D11").Select()
